$d = $word.ActiveDocument

# Pull the whole package as a single OOXML string so we can surgically patch
# both word/document.xml and word/styles.xml in one shot, then feed the
# modified package back in through Range.InsertXML (applying to the full
# document content addresses every part in the package, unlike InsertXML on
# a sub-range which only touches word/document.xml).
$pkg = $d.Content.XML()

# ---------------------------------------------------------------------
# 1) Title paragraph: "Lazatech Courses" -> spell-checked split runs
# ---------------------------------------------------------------------
$oldTitleRun = '<w:r w:rsidRPr="00396DF7"><w:t>Lazatech Courses</w:t></w:r>'
$newTitleRun = '<w:proofErr w:type="spellStart"/><w:r><w:t>Lazatech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Courses</w:t></w:r>'
if ($pkg.IndexOf($oldTitleRun) -lt 0) { throw "title run not found" }
$pkg = $pkg.Replace($oldTitleRun, $newTitleRun)

# ---------------------------------------------------------------------
# 2) Welcome paragraph: split out "Lazatech" for spell-check marks, and
#    append the new "Course 1" block after the existing sentence.
# ---------------------------------------------------------------------
$oldWelcomeRun = '<w:r w:rsidRPr="00D1596F"><w:t xml:space="preserve">               Welcome to the Lazatech Educate website. All of the available courses are mentioned here. You are free to enroll in any course you like. You are also free to select the lessons you desire to complete. You may learn at your own speed, and you can track your progress and earnings using the navigation bar above.</w:t></w:r>'

$newWelcomeBlock = '<w:r><w:t xml:space="preserve">               Welcome to the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Lazatech</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> Educate website. All of the available courses are mentioned here. You are free to enroll in any course you like. You are also free to select the lessons you desire to complete. You may learn at your own speed, and you can track your progress and earnings using the navigation bar above.</w:t></w:r>' +
  '<w:r><w:br/></w:r>' +
  '<w:r><w:br/></w:r>' +
  '<w:r><w:rPr><w:rStyle w:val="Heading2Char"/></w:rPr><w:t>Course 1: Getting to know LAZATECH EDUCATE</w:t></w:r>' +
  '<w:r><w:rPr><w:rStyle w:val="Heading2Char"/></w:rPr><w:t xml:space="preserve"> AND ITS SOCIAL MEDIA COUNTERPART</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r>' +
  '<w:r><w:t xml:space="preserve">Complete this course if you want to get the most out of </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Lazatech</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> Educate and our ecosystem, since it covers the project as a whole as well as the suite of utilities we are introducing to the XRPL.</w:t></w:r>' +
  '<w:r><w:br/></w:r>'

if ($pkg.IndexOf($oldWelcomeRun) -lt 0) { throw "welcome run not found" }
$pkg = $pkg.Replace($oldWelcomeRun, $newWelcomeBlock)

# ---------------------------------------------------------------------
# 3) styles.xml: add the "Heading2" paragraph style (right after Heading1,
#    mirroring where Word places newly-used built-in styles) and the
#    linked "Heading2Char" character style (right after Heading1Char).
# ---------------------------------------------------------------------
$heading1Style = '<w:style w:type="paragraph" w:styleId="Heading1"><w:name w:val="heading 1"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00BB0B5F"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="240" w:after="0"/><w:outlineLvl w:val="0"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:style>'

$heading2Style = '<w:style w:type="paragraph" w:styleId="Heading2"><w:name w:val="heading 2"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading2Char"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="00372EBA"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="40" w:after="0"/><w:outlineLvl w:val="1"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:style>'

if ($pkg.IndexOf($heading1Style) -lt 0) { throw "Heading1 style not found" }
$pkg = $pkg.Replace($heading1Style, $heading1Style + $heading2Style)

$heading1CharStyle = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading1Char"><w:name w:val="Heading 1 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="00BB0B5F"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:style>'

$heading2CharStyle = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading2Char"><w:name w:val="Heading 2 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading2"/><w:uiPriority w:val="9"/><w:rsid w:val="00372EBA"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:style>'

if ($pkg.IndexOf($heading1CharStyle) -lt 0) { throw "Heading1Char style not found" }
$pkg = $pkg.Replace($heading1CharStyle, $heading1CharStyle + $heading2CharStyle)

# ---------------------------------------------------------------------
# Write everything back in a single pass.
# ---------------------------------------------------------------------
$d.Content.InsertXML($pkg)

Write-Output "done"
